$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the first data row (row 2), shifting all existing
# data rows down by one (old row 2 -> row 3, ... old row 48 -> row 49).
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting (bold/border) from the header row
# above it; clear that so the new row matches the plain style used by the
# rest of the data rows.
$ws.Rows.Item(2).ClearFormats()

# Populate the new row 2 with the new weekly price entry.
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C2").Value = "Los Lagos"
$ws.Range("D2").Value = 44956
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100101
$ws.Range("H2").Value = "Berries"
$ws.Range("I2").Value = 100101001
$ws.Range("J2").Value = "Arándano (blue)"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 2000
$ws.Range("O2").Value = 2200
$ws.Range("P2").Value = 2100
$ws.Range("Q2").Value = "$/bandeja 2 kilos"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1050
$ws.Range("T2").Value = 2

# Match the date-formatted style used by the rest of the "Fecha" column.
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat
